$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6: Time changes from "1000-1230" to "1000-1300", Hours from 2.5 to 3
$ws.Range("B6").Value = "1000-1300"
$ws.Range("C6").Value = 3

# Update selection to match the new active cell
$ws.Range("B5").Select()
